# Append order-inventory rows 18-29 (order numbers 17-26) to Sheet1.
# Existing sheet data runs through row 17 (A1:F17); this appends through F29.
# Every value in this sheet is stored as literal text (inlineStr in the
# source file), including numeric-looking and "$"-prefixed values. Prefixing
# each value with a leading apostrophe forces Excel's quote-prefix / literal
# text entry, so numeric-looking strings ("17", "6000", "$2000", ...) land as
# text instead of being auto-coerced to numbers/currency; for values that are
# already non-numeric text (e.g. "mouse", "Hanna") the apostrophe is simply
# the standard "treat as text" marker and does not appear in the stored value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("17", "2320232342342", "mouse",  "10",    "Hanna", "Ian"),
    @("17", "2320232342342", "mouse",  "10",    "Hanna", "Ian"),
    @("18", "2320232342342", "cat",    "100",   "Mike",  "Ian"),
    @("18", "2320232342342", "dog",    "100",   "Mike",  "Ian"),
    @("19", "2320232342342", "box",    "1000",  "Greg",  "Ian"),
    @("20", "2320232342342", "circle", "$2000", "Greg",  "Ian"),
    @("21", "2320232342342", "square", "3000",  "Greg",  "Ian"),
    @("22", "2320232342342", "square", "5000",  "Greg",  "Ian"),
    @("23", "2320232342342", "square", "6000",  "greg",  "ian"),
    @("24", "2320232342342", "square", "6000",  "greg",  "ian"),
    @("25", "2320232342342", "square", "6000",  "greg",  "ian"),
    @("26", "2320232342342", "square", "6000",  "greg",  "ian")
)

$startRow = 18
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = "'" + $values[$c - 1]
    }
}
